# Updates the "cryptos" list with refreshed prices / 1h volume percentages,
# and swaps the Filecoin / HuobiToken rows (rows 33 and 34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = "27.226.31"
$dCell.Style = "Normal"
$ws.Range("E2").Value = "  -3.48%  "

# --- Row 3 ---
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = "1.810.67"
$dCell.Style = "Normal"
$ws.Range("E3").Value = "  -3.54%  "

# --- Row 4 ---
$ws.Range("E4").Value = "  +0.00%  "

# --- Row 5 ---
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = "310.28"
$dCell.Style = "Normal"

# --- Row 6 ---
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E6").Value = "  +0.03%  "

# --- Row 7 ---
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4203"
$dCell.Style = "Normal"

# --- Row 8 ---
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = "0.3558"
$dCell.Style = "Normal"
$ws.Range("E8").Value = "  -4.04%  "

# --- Row 9 ---
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = "0.07109"
$dCell.Style = "Normal"
$ws.Range("E9").Value = "  -4.17%  "

# --- Row 10 ---
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = "0.8501"
$dCell.Style = "Normal"
$ws.Range("E10").Value = "  -3.67%  "

# --- Row 11 ---
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = "20.21"
$dCell.Style = "Normal"
$ws.Range("E11").Value = "  -4.30%  "

# --- Row 12 ---
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = "1.794.81"
$dCell.Style = "Normal"
$ws.Range("E12").Value = "  -5.79%  "

# --- Row 13 ---
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = "5.308"
$dCell.Style = "Normal"
$ws.Range("E13").Value = "  -3.37%  "

# --- Row 14 ---
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = "6.389"
$dCell.Style = "Normal"
$ws.Range("E14").Value = "  -3.59%  "

# --- Row 15 ---
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06860"
$dCell.Style = "Normal"
$ws.Range("E15").Value = "  -1.85%  "

# --- Row 16 ---
$ws.Range("E16").Value = "  +0.04%  "

# --- Row 17 ---
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = "81.17"
$dCell.Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "

# --- Row 18 ---
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = "0.000008730"
$dCell.Style = "Normal"
$ws.Range("E18").Value = "  -4.51%  "

# --- Row 19 ---
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "

# --- Row 20 ---
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = "15.14"
$dCell.Style = "Normal"
$ws.Range("E20").Value = "  -3.13%  "

# --- Row 21 ---
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = "27.137.84"
$dCell.Style = "Normal"
$ws.Range("E21").Value = "  -3.89%  "

# --- Row 22 ---
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = "5.123"
$dCell.Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

# --- Row 23 ---
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = "10.85"
$dCell.Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "

# --- Row 24 ---
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = "2.002.93"
$dCell.Style = "Normal"
$ws.Range("E24").Value = "  -5.72%  "

# --- Row 25 ---
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = "1.974"
$dCell.Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "

# --- Row 26 ---
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = "153.79"
$dCell.Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "

# --- Row 27 ---
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = "18.15"
$dCell.Style = "Normal"
$ws.Range("E27").Value = "  -3.25%  "

# --- Row 28 ---
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = "5.062"
$dCell.Style = "Normal"
$ws.Range("E28").Value = "  -6.77%  "

# --- Row 29 ---
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = "113.44"
$dCell.Style = "Normal"
$ws.Range("E29").Value = "  -3.56%  "

# --- Row 30 ---
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = "1.701"
$dCell.Style = "Normal"
$ws.Range("E30").Value = "  -9.14%  "

# --- Row 31 ---
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = "0.08900"
$dCell.Style = "Normal"
$ws.Range("E31").Value = "  -0.85%  "

# --- Row 32 ---
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = "0.7417"
$dCell.Style = "Normal"
$ws.Range("E32").Value = "  -6.51%  "

# --- Row 33 ---
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = "2.918"
$dCell.Style = "Normal"
$ws.Range("E33").Value = "  -1.09%  "

# --- Row 34 ---
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = "4.445"
$dCell.Style = "Normal"
$ws.Range("E34").Value = "  -5.83%  "

# --- Row 35 ---
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = "1.104"
$dCell.Style = "Normal"
$ws.Range("E35").Value = "  -7.12%  "

# --- Row 36 ---
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = "1.002"
$dCell.Style = "Normal"
$ws.Range("E36").Value = "  +0.10%  "

# --- Row 37 ---
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = "1.068"
$dCell.Style = "Normal"
$ws.Range("E37").Value = "  -5.42%  "

# --- Row 38 ---
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = "0.05199"
$dCell.Style = "Normal"
$ws.Range("E38").Value = "  -4.85%  "

# --- Row 39 ---
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = "0.01903"
$dCell.Style = "Normal"
$ws.Range("E39").Value = "  -2.93%  "

# --- Row 40 ---
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = "0.1636"
$dCell.Style = "Normal"
$ws.Range("E40").Value = "  -3.34%  "

# --- Row 41 ---
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = "2.702"
$dCell.Style = "Normal"
$ws.Range("E41").Value = "  -6.39%  "

# --- Row 42 ---
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4964"
$dCell.Style = "Normal"
$ws.Range("E42").Value = "  -4.04%  "

# --- Row 43 ---
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = "6.273"
$dCell.Style = "Normal"
$ws.Range("E43").Value = "  -8.79%  "

# --- Row 44 ---
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = "8.157"
$dCell.Style = "Normal"
$ws.Range("E44").Value = "  -5.35%  "

# --- Row 45 ---
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = "105.14"
$dCell.Style = "Normal"
$ws.Range("E45").Value = "  -0.57%  "

# --- Row 46 ---
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = "10.17"
$dCell.Style = "Normal"
$ws.Range("E46").Value = "  -4.38%  "

# --- Row 47 ---
$ws.Range("E47").Value = "  +0.09%  "

# --- Row 48 ---
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = "0.06383"
$dCell.Style = "Normal"
$ws.Range("E48").Value = "  -3.14%  "

# --- Row 49 ---
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = "0.4563"
$dCell.Style = "Normal"
$ws.Range("E49").Value = "  -4.35%  "

# --- Row 50 ---
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = "1.592"
$dCell.Style = "Normal"
$ws.Range("E50").Value = "  -3.71%  "

# --- Row 51 ---
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = "62.92"
$dCell.Style = "Normal"
$ws.Range("E51").Value = "  -3.94%  "
